# Update market-price derived figures (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, as produced by the scheduled
# price-refresh runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 396.77777
$ws.Range("I42").Value = 346.375
$ws.Range("K42").Value = 1039.125
$ws.Range("M42").Value = -809.125

$ws.Range("H51").Value = 4346.75
$ws.Range("J51").Value = 4346.75
$ws.Range("L51").Value = 4346.75
$ws.Range("N51").Value = -5314.75

$ws.Range("H69").Value = 3973
$ws.Range("J69").Value = 3979.5
$ws.Range("L69").Value = 11938.5
$ws.Range("N69").Value = -13686.5

$ws.Range("H72").Value = 3973
$ws.Range("J72").Value = 3979.5
$ws.Range("L72").Value = 35815.5
$ws.Range("N72").Value = -44551.5

$ws.Range("H74").Value = 6788.625
$ws.Range("I74").Value = 6788.625
$ws.Range("K74").Value = 6788.625
$ws.Range("M74").Value = -5852.625

$ws.Range("H77").Value = 6788.625
$ws.Range("I77").Value = 6788.625
$ws.Range("K77").Value = 33943.125
$ws.Range("M77").Value = -29263.125

$ws.Range("H98").Value = 2898.9048
$ws.Range("I98").Value = 2530.3125
$ws.Range("J98").Value = 4078.4
$ws.Range("K98").Value = 2530.3125
$ws.Range("L98").Value = 4078.4
$ws.Range("M98").Value = -1032.3125
$ws.Range("N98").Value = -7074.4

$ws.Range("H122").Value = 2898.9048
$ws.Range("I122").Value = 2530.3125
$ws.Range("J122").Value = 4078.4
$ws.Range("K122").Value = 7590.9375
$ws.Range("L122").Value = 12235.2
$ws.Range("M122").Value = -5140.9375
$ws.Range("N122").Value = -17135.2

$ws.Range("H125").Value = 1554.8
$ws.Range("I125").Value = 1776.5
$ws.Range("J125").Value = 668
$ws.Range("K125").Value = 15988.5
$ws.Range("L125").Value = 6012
$ws.Range("M125").Value = -13528.5
$ws.Range("N125").Value = -10932

$ws.Range("H137").Value = 1285.9032
$ws.Range("I137").Value = 1193.9166
$ws.Range("J137").Value = 1601.2858
$ws.Range("K137").Value = 3581.7498
$ws.Range("L137").Value = 4803.857400000001
$ws.Range("M137").Value = -1031.7498
$ws.Range("N137").Value = -9903.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 977.6
$ws.Range("I25").Value = 977.6
$ws.Range("K25").Value = 977.6
$ws.Range("M25").Value = -575.6

$ws.Range("H61").Value = 1234.5294
$ws.Range("I61").Value = 963.4286
$ws.Range("J61").Value = 2499.6667
$ws.Range("K61").Value = 963.4286
$ws.Range("L61").Value = 2499.6667
$ws.Range("M61").Value = -751.4286
$ws.Range("N61").Value = -2923.6667

$ws.Range("H110").Value = 2241.8572
$ws.Range("I110").Value = 1126.6666
$ws.Range("K110").Value = 1126.6666
$ws.Range("M110").Value = 918.3334

$ws.Range("H136").Value = 1234.5294
$ws.Range("I136").Value = 963.4286
$ws.Range("J136").Value = 2499.6667
$ws.Range("K136").Value = 2890.2858
$ws.Range("L136").Value = 7499.000100000001
$ws.Range("M136").Value = -340.2857999999997
$ws.Range("N136").Value = -12599.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 140.8
$ws.Range("J5").Value = 100
$ws.Range("L5").Value = 100
$ws.Range("N5").Value = -326

$ws.Range("H43").Value = 237342
$ws.Range("J43").Value = 237342
$ws.Range("L43").Value = 237342
$ws.Range("N43").Value = -237704

$ws.Range("H61").Value = 17000
$ws.Range("J61").Value = 17000
$ws.Range("L61").Value = 17000
$ws.Range("N61").Value = -17626

$ws.Range("H107").Value = 1614.8572
$ws.Range("I107").Value = 1012.25
$ws.Range("K107").Value = 1012.25
$ws.Range("M107").Value = 907.75

$ws.Range("H134").Value = 7975.5557
$ws.Range("I134").Value = 968.6429000000001
$ws.Range("J134").Value = 32499.75
$ws.Range("K134").Value = 2905.9287
$ws.Range("L134").Value = 97499.25
$ws.Range("M134").Value = -370.9287000000004
$ws.Range("N134").Value = -102569.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H57").Value = 16000
$ws.Range("J57").Value = 16000
$ws.Range("L57").Value = 16000
$ws.Range("N57").Value = -17120

$ws.Range("H122").Value = 2800.4
$ws.Range("I122").Value = 1004
$ws.Range("J122").Value = 5495
$ws.Range("K122").Value = 3012
$ws.Range("L122").Value = 16485
$ws.Range("M122").Value = -562
$ws.Range("N122").Value = -21385

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1524.25
$ws.Range("I5").Value = 1524.25
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4572.75
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4460.75
$ws.Range("N5").ClearContents()

$ws.Range("H135").Value = 1524.25
$ws.Range("I135").Value = 1524.25
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 13718.25
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -11183.25
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16000
$ws.Range("I15").Value = 10000
$ws.Range("K15").Value = 10000
$ws.Range("M15").Value = -9712

$ws.Range("H70").Value = 26476778
$ws.Range("I70").Value = 20839100
$ws.Range("K70").Value = 20839100
$ws.Range("M70").Value = -20838830

$ws.Range("H73").Value = 26476778
$ws.Range("I73").Value = 20839100
$ws.Range("K73").Value = 20839100
$ws.Range("M73").Value = -20838164

$ws.Range("H80").Value = 4702.4287
$ws.Range("I80").Value = 3435
$ws.Range("J80").Value = 5653
$ws.Range("K80").Value = 3435
$ws.Range("L80").Value = 5653
$ws.Range("M80").Value = -2437
$ws.Range("N80").Value = -7649

$ws.Range("H81").Value = 16000
$ws.Range("I81").Value = 10000
$ws.Range("K81").Value = 10000
$ws.Range("M81").Value = -9002

$ws.Range("H83").Value = 4702.4287
$ws.Range("I83").Value = 3435
$ws.Range("J83").Value = 5653
$ws.Range("K83").Value = 17175
$ws.Range("L83").Value = 28265
$ws.Range("M83").Value = -12183
$ws.Range("N83").Value = -38249

$ws.Range("H84").Value = 16000
$ws.Range("I84").Value = 10000
$ws.Range("K84").Value = 30000
$ws.Range("M84").Value = -25008

$ws.Range("H113").Value = 3233.5
$ws.Range("I113").Value = 2207.3333
$ws.Range("J113").Value = 3673.2856
$ws.Range("K113").Value = 2207.3333
$ws.Range("L113").Value = 3673.2856
$ws.Range("M113").Value = -37.33329999999978
$ws.Range("N113").Value = -8013.2856

$ws.Range("H122").Value = 151685.7
$ws.Range("I122").Value = 2078.4285
$ws.Range("J122").Value = 500769.34
$ws.Range("K122").Value = 6235.2855
$ws.Range("L122").Value = 1502308.02
$ws.Range("M122").Value = -3785.2855
$ws.Range("N122").Value = -1507208.02

$ws.Range("H132").Value = 2826.9
$ws.Range("I132").Value = 2496.889
$ws.Range("K132").Value = 7490.667
$ws.Range("M132").Value = -4960.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 1913.5
$ws.Range("I68").Value = 1668.3334
$ws.Range("K68").Value = 1668.3334
$ws.Range("M68").Value = -919.3334

$ws.Range("H71").Value = 1913.5
$ws.Range("I71").Value = 1668.3334
$ws.Range("K71").Value = 8341.666999999999
$ws.Range("M71").Value = -4597.666999999999

$ws.Range("H132").Value = 20066.926
$ws.Range("I132").Value = 1303.3667
$ws.Range("J132").Value = 43521.375
$ws.Range("K132").Value = 3910.1001
$ws.Range("L132").Value = 130564.125
$ws.Range("M132").Value = -1380.1001
$ws.Range("N132").Value = -135624.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 54250
$ws.Range("J127").Value = 54250
$ws.Range("L127").Value = 54250
$ws.Range("N127").Value = -64170

$ws.Range("H132").Value = 3085.8333
$ws.Range("I132").Value = 2597.8235
$ws.Range("J132").Value = 4271
$ws.Range("K132").Value = 7793.470499999999
$ws.Range("L132").Value = 12813
$ws.Range("M132").Value = -5263.470499999999
$ws.Range("N132").Value = -17873
Write-Host "All updates applied"
